$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SPN")
$ws1.Activate()
$ws1.Range("K26").Select()
$excel.Selection.Resize(0,2).Select()
Write-Host $excel.ActiveCell.Address()
Write-Host $excel.Selection.Address()
